# calc_shear.xlsx -- span 2 ltr/rtl test data fix-up
#
# - swap the Pr (C3) / Pl (C4) values that had been entered backwards
# - flip the "direction" flag (C5) from ltr to rtl
# - move the active selection off C7 onto C4 (where the edit actually is)
# - give the Ve result cell (C7) a numeric "0.0000" display format
# - widen the sheet-tab area of the window splitter (tabRatio)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pr / Pl were transposed -- put them back the right way round.
$ws.Range("C3").Value = 360
$ws.Range("C4").Value = 248

# direction flag: ltr -> rtl
$ws.Range("C5").Value = "rtl"

# Result cell now gets an explicit 4-decimal numeric format.
$ws.Range("C7").NumberFormat = "0.0000"

# Move the cursor/selection to C4 (the cell that was actually being worked on).
$ws.Range("C4").Select()

# Widen the tab-bar / horizontal-scrollbar split in the window (cosmetic).
$excel.ActiveWindow.TabRatio = 472
